# edit.ps1 -- apply the "Added Ryan Keller and Andy Ko to opening slides." commit
# to the currently open presentation ($ppt.ActivePresentation).

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 1 ("INFO 461: Cooperative Software Development") - Subtitle box:
# add a "TA:" line (Ryan Keller) and an "LA:" line (Prof. Andy Ko) below
# the existing "Instructor:" line, and widen the box to fit the new text.
# -----------------------------------------------------------------------
$s1  = $p.Slides.Item(1)
$sub = $s1.Shapes.Item(2)

$sub.TextFrame.TextRange.Text = "Instructor: Dr. Andrew Begel, Microsoft Research abegel@uw.edu" + "`r" + "TA: Ryan Keller, UW rykeller@uw.edu" + "`r" + "LA: Prof. Andy Ko, UW ajko@uw.edu"

$subTr = $sub.TextFrame.TextRange
$subTr.Characters(50, 13).Font.Italic = $true    # abegel@uw.edu
$subTr.Characters(84, 15).Font.Italic = $true    # rykeller@uw.edu
$subTr.Characters(122, 11).Font.Italic = $true   # ajko@uw.edu

# Widen the subtitle placeholder so the longer lines fit.
$sub.Width = 874.64

# -----------------------------------------------------------------------
# Slide 9 ("Today's Class") - 4th agenda bullet:
# "Effort Estimation" -> "SMART Commitments"
# -----------------------------------------------------------------------
$s9  = $p.Slides.Item(9)
$agenda = $s9.Shapes.Item(2)
$agenda.TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "SMART Commitments"

# -----------------------------------------------------------------------
# Slide 10 ("How this Course Works") - feedback bullet:
# "Let me know! I want your feedback..." -> "Let us know! We want your feedback..."
# -----------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$body10 = $s10.Shapes.Item(2)
$body10.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Let us know! We want your feedback to improve the course."
